$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 41; this shifts existing rows 41-97 down to 42-98
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new record's values
$ws.Cells.Item(41, 1).Value = 4
$ws.Cells.Item(41, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(41, 3).Value = "Los Lagos"
$ws.Cells.Item(41, 4).Value = 44803
$ws.Cells.Item(41, 5).Value = 10
$ws.Cells.Item(41, 6).Value = 100112026
$ws.Cells.Item(41, 7).Value = "Haba"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 80
$ws.Cells.Item(41, 11).Value = 16000
$ws.Cells.Item(41, 12).Value = 16000
$ws.Cells.Item(41, 13).Value = 16000
$ws.Cells.Item(41, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(41, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(41, 16).Value = 640
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
